# The "Total" row on the Gesamtergebnis sheet previously reported "N/A" for
# the balance columns (Startguthaben / Endsaldo) because summing balances
# across platforms wasn't supported. The writer now also totals those two
# columns, so replace the "N/A" placeholders with the actual column sums.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gesamtergebnis")
$ws.Activate()

# Total = sum of the per-platform rows (row 2: Bondora, row 3: DoFinance)
$ws.Range("C4").Value = $ws.Range("C2").Value2 + $ws.Range("C3").Value2
$ws.Range("D4").Value = $ws.Range("D2").Value2 + $ws.Range("D3").Value2

# Leave the selection on the newly-updated totals, as a user reviewing the
# change would.
$ws.Range("C5").Select()
